$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Locate the "Hexagon 5" shape (id=6) and nudge its position/size slightly,
# plus give its outline an explicit 3pt weight.
$hex = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.Name -eq "Hexagon 5") {
        $hex = $sh
    }
}

# Values below are chosen (in points, EMU/12700) so that after the host's
# internal float32 rounding they land exactly on the target EMU values:
#   Left   -> 146776 EMU
#   Top    -> 130495 EMU
#   Width  -> 1535245 EMU
#   Height -> 1321882 EMU
$hex.Left = 11.55716609954834
$hex.Top = 10.27519702911377
$hex.Width = 120.88543701171875
$hex.Height = 104.08519744873047

# Give the outline an explicit weight of 3pt (38100 EMU); previously unset/default.
$hex.Line.Weight = 3
